# Updated symbol list on Tue Dec 13 04:49:42 UTC 2022 with GitHub Actions
#
# Applies the per-row "Price" (column D) refreshes plus the BKEXToken/CEJI
# row swap (rows 42-43) to the crypto price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Cell, $Text)
    # Force text storage so numeric-looking strings (e.g. "0.1640",
    # "0.004100") keep their exact formatting instead of being coerced into
    # a Double (which would drop significant trailing/leading zeros).
    $rng = $Worksheet.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
}

# ---- Column D ("Price") refreshes ----
Set-TextValue $ws "D2"  "268.37"
Set-TextValue $ws "D4"  "6.244"
Set-TextValue $ws "D5"  "0.06210"
Set-TextValue $ws "D6"  "3.569"
Set-TextValue $ws "D7"  "6.544"
Set-TextValue $ws "D8"  "1.393"
Set-TextValue $ws "D9"  "0.8269"
Set-TextValue $ws "D10" "0.1640"
Set-TextValue $ws "D11" "0.08277"
Set-TextValue $ws "D12" "0.03569"
Set-TextValue $ws "D13" "0.03191"
Set-TextValue $ws "D14" "0.09202"
Set-TextValue $ws "D15" "3.769"
Set-TextValue $ws "D17" "0.04670"
Set-TextValue $ws "D18" "0.006417"
Set-TextValue $ws "D19" "0.006206"
Set-TextValue $ws "D23" "2.270"
Set-TextValue $ws "D24" "0.01366"
Set-TextValue $ws "D28" "0.0002713"
Set-TextValue $ws "D40" "0.04718"
Set-TextValue $ws "D41" "0.006972"

# ---- Rows 42/43: BKEXToken and CEJI swap places ----
# Row 42 was BKEXToken, now becomes CEJI.
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.004100"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 was CEJI, now becomes BKEXToken.
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1122"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# ---- Remaining column D ("Price") refreshes ----
Set-TextValue $ws "D44" "0.01160"
Set-TextValue $ws "D45" "0.00006265"
Set-TextValue $ws "D46" "0.0009899"
Set-TextValue $ws "D48" "0.8026"
Set-TextValue $ws "D49" "0.002335"
